$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (columns C/D/E) ---

# Row 8 (extr1): C 5->14, D 12->11, E False->True
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 (extr2): C 5->16, E False->True
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# Row 10 (extr3): C 10->5, D 11->12
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11 (extr4): C 7->5, D 8->9
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12 (extr5): C 9->10
$ws.Range("C12").Value = 10

# Row 13 (extr6): D 11->8
$ws.Range("D13").Value = 8

# Row 14 (extr7): C 5->9, D 7->11
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# Row 15 (extr8): C 8->7, D 5->11, E True->False
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- Append two new rows (16 and 17), copying formatting from row 15 ---

$ws.Range("A15:E15").Copy()
$ws.Range("A16:E17").PasteSpecial(-4122)

# Row 16: line7
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# Row 17: line8
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
